# Automated hashcode update
# Updates the hashcode values (column B) for a set of rows in the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = "6dcfb0297c0f89cab5a3ed3a0f747851"
    11  = "46f4b710c85b0f0cc5f14b22ef36eacc"
    29  = "b0baac0fd6228bbfd725ed8e9176549f"
    32  = "9331c4bf41bae51046f4585b269baaa2"
    121 = "aa7b52c533406a158a3af883053ad5d9"
    123 = "46ed7460e568b32a9134c96746c4b30b"
    126 = "e76deb34b2d2c58dc4cddf2e5ec96a06"
    175 = "5711ec02e10e6127e2cabd7d619d26a4"
    213 = "6ff8568f98d8581660b65262094e9bd3"
    228 = "8a5eff7babfb9adfead2564f6f30314e"
    351 = "f08056f7dc8bdee4de507972d9f93567"
    354 = "03f96168b2e3e4ebad6414a1b23f1845"
    397 = "af0d8651b105f6374a0fcd8e37c58b94"
    402 = "77ad9782706a0a7bd9d04bc9097ffe21"
    461 = "b2ed4f7be4ba8ef8c2eca34e9c152743"
    513 = "922820cb546d4143611e0ac0c6cb3e5c"
    521 = "4725d2dc189712fda585ce4142710523"
    572 = "f84b2dc7a60816718e52ec71e638e166"
    629 = "db690ac0b9e5d4c7fcabcea242f678b2"
    649 = "7c9e1afc19068480855439254cec2b6b"
    655 = "b8fe97b8a9e6a16a0b340ab282597ab0"
    715 = "848fb2e1daaef48a2baee487df02ad08"
    781 = "3530be274c9da14179c1054bb965cea0"
    788 = "a2571fcd86de04f7e6f45ed90ca857b4"
    862 = "84cca520ec83e9af891d6cfa4ab2a952"
    874 = "d878f735a89572d2273c1e98708e28dd"
    896 = "dc5ab44aaf01eeca4909629fce968836"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
